$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber colaborador_id / treinamentos_id columns (A2:B6 and A7:B7)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 6

# Row 7 (C7/D7/E7) was carrying a stray font/number-format combo; bring it
# back in line with the rest of the table by pulling the formatting used
# on the corresponding columns elsewhere, then set D7's value to the
# "Completo" status text.
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "Completo"

$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Move the active selection to B5
$ws.Range("B5").Select() | Out-Null
